# Weekly fruit/vegetable price update for "Fruta, Vega Central Mapocho de Santiago - Coco".
# Each data row (2-41) is refreshed with a new report date and the associated
# volume/min/max/weighted-avg prices, unit price per kg, quality and origin,
# matching this week's consolidated logica_diaria extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Row, Fecha(D), Calidad(L), Volumen(M), PrecioMin(N), PrecioMax(O), PrecioProm(P), Origen(R), PrecioKg(S)
$rows = @(
    @(2, 44435, "Primera", 60, 25000, 25000, 25000, "Perú", 1250),
    @(3, 44410, "Primera", 40, 25000, 25000, 25000, "Perú", 1250),
    @(4, 44166, "Primera", 120, 28000, 28000, 28000, "Perú", 1400),
    @(5, 44438, "Primera", 25, 21000, 21000, 21000, "Perú", 1050),
    @(6, 44277, "Primera", 60, 24000, 24000, 24000, "Perú", 1200),
    @(7, 44431, "Primera", 60, 25000, 25000, 25000, "Perú", 1250),
    @(8, 44522, "Primera", 25, 30000, 30000, 30000, "Perú", 1500),
    @(9, 44396, "Primera", 45, 22000, 22000, 22000, "Perú", 1100),
    @(10, 44473, "Primera", 40, 24000, 24000, 24000, "Perú", 1200),
    @(11, 44363, "Primera", 150, 21000, 22000, 21500, "Perú", 1075),
    @(12, 44417, "Primera", 30, 24000, 24000, 24000, "Perú", 1200),
    @(13, 44830, "Primera", 200, 30000, 30000, 30000, "Perú", 1500),
    @(14, 44442, "Primera", 30, 22000, 22000, 22000, "Perú", 1100),
    @(15, 44372, "Primera", 60, 20000, 21000, 20667, "Perú", 1033),
    @(16, 44298, "Primera", 240, 19000, 20000, 19500, "Perú", 975),
    @(17, 44284, "Primera", 40, 23000, 23000, 23000, "Perú", 1150),
    @(18, 44305, "Primera", 40, 24000, 24000, 24000, "Perú", 1200),
    @(19, 44299, "Primera", 150, 19000, 20000, 19500, "Perú", 975),
    @(20, 44263, "Segunda", 150, 15000, 15000, 15000, "Perú", 750),
    @(21, 44354, "Primera", 150, 21000, 22000, 21500, "Perú", 1075),
    @(22, 44445, "Primera", 35, 20000, 20000, 20000, "Perú", 1000),
    @(23, 44302, "Primera", 100, 19000, 20000, 19500, "Perú", 975),
    @(24, 44365, "Primera", 150, 20000, 21000, 20500, "Perú", 1025),
    @(25, 44312, "Primera", 50, 22000, 22000, 22000, "Perú", 1100),
    @(26, 44357, "Primera", 200, 20000, 21000, 20500, "Perú", 1025),
    @(27, 44300, "Primera", 150, 19000, 20000, 19500, "Perú", 975),
    @(28, 44382, "Primera", 200, 19000, 20000, 19500, "Perú", 975),
    @(29, 44326, "Primera", 40, 22000, 22000, 22000, "Perú", 1100),
    @(30, 44760, "Primera", 300, 24000, 25000, 24500, "Perú", 1225),
    @(31, 44529, "Primera", 34, 28000, 28000, 28000, "Perú", 1400),
    @(32, 44452, "Primera", 35, 21000, 22000, 21429, "Perú", 1071),
    @(33, 44270, "Primera", 50, 24000, 24000, 24000, "Perú", 1200),
    @(34, 44165, "Primera", 300, 27000, 28000, 27500, "Perú", 1375),
    @(35, 44424, "Primera", 70, 24000, 25000, 24429, "Perú", 1221),
    @(36, 44355, "Primera", 200, 20000, 21000, 20500, "Ecuador", 1025),
    @(37, 44356, "Primera", 100, 20000, 21000, 20500, "Perú", 1025),
    @(38, 44350, "Primera", 90, 21000, 22000, 21556, "Perú", 1078),
    @(39, 44333, "Primera", 30, 22000, 22000, 22000, "Perú", 1100),
    @(40, 44620, "Primera", 60, 22000, 22000, 22000, "Perú", 1100),
    @(41, 44613, "Primera", 60, 30000, 30000, 30000, "Perú", 1100)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("D$r").Value = $row[1]
    $ws.Range("L$r").Value = $row[2]
    $ws.Range("M$r").Value = $row[3]
    $ws.Range("N$r").Value = $row[4]
    $ws.Range("O$r").Value = $row[5]
    $ws.Range("P$r").Value = $row[6]
    $ws.Range("R$r").Value = $row[7]
    $ws.Range("S$r").Value = $row[8]
}
